# Add "Norway" and "Poland" market test-data sheets, modeled on the
# existing "Croatia" sheet (same layout/column widths/no extra row
# heights as the sheets that will host them), placed after "Hungary".

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("Croatia")

# --- Norway --------------------------------------------------------
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$norway = $wb.Worksheets.Item($wb.Worksheets.Count)
$norway.Name = "Norway"
$norway.Range("B4").Value = "NGC-2931/T3079"
$norway.Range("B2").Value = "Norway Market"

# --- Poland ----------------------------------------------------------
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$poland = $wb.Worksheets.Item($wb.Worksheets.Count)
$poland.Name = "Poland"
$poland.Range("B4").Value = "NGC-2920/T3108"
$poland.Range("B2").Value = "Poland Market"

# Norway ends up the active tab (matches the captured edit).
$norway.Activate()
